# "adicionando alterações finais da sprint 2 entregue 19.09"
#
# 1) Duplicate the (still-original) SPRINT 2 block into a new SPRINT 3
#    block (rows 21-34), before SPRINT 2 itself gets edited, so the
#    copy starts from a clean, identically-formatted template.
# 2) Update the SPRINT 2 block (rows 7-19): delivery date, status, and
#    responsible/"next sprint" annotations in column C.
# 3) Overwrite the freshly copied SPRINT 3 block with its own tasks /
#    responsibles / status.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Clone SPRINT 2's layout/formatting into the SPRINT 3 rows
# ---------------------------------------------------------------------

# Title + header rows + 8 "plain" task rows (row 19's footer-style
# bottom border is intentionally excluded here; it's rebuilt below).
$ws.Range("A7:C18").Copy($ws.Range("A21"))

# New 9th task row, re-using the plain (non-footer) row style.
$ws.Range("B32:C32").Copy($ws.Range("B33"))

# Final task row, re-using the footer-style (bottom border) row.
$ws.Range("A19:C19").Copy($ws.Range("A34"))

# ---------------------------------------------------------------------
# 2) SPRINT 2 updates
# ---------------------------------------------------------------------

# Delivery date moved from 09/17/2019 to 09/19/2019
$ws.Range("B9").Value2 = 43727

# Status updated - now delivered, so restyle from the red "in progress"
# look to the bold "Entregue!" look already used up in B5.
$ws.Range("B10").Value2 = "Entregue! (com algumas coisas faltantes)"
$ws.Range("B5").Copy()
$ws.Range("B10").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Column C "Responsável" notes for each page task
$ws.Range("C12").Value2 = "Laura - OK"
$ws.Range("C13").Value2 = "Stela - OK"
$ws.Range("C14").Value2 = "Higor - OK"
$ws.Range("C15").Value2 = "Sabrine - OK"
$ws.Range("C16").Value2 = "Higor - OK"

# Task names refined with category prefixes
$ws.Range("B17").Value2 = "IMAGENS: Seleção de imagens para todas as páginas"
$ws.Range("B18").Value2 = "TEXTOS E TÍTULOS: Definição dos textos para todas as páginas"
$ws.Range("B19").Value2 = "LOGO: Elaboração do Logo"

# These three tasks slipped to the next sprint - flag them in red italics
$ws.Range("C17").Value2 = "Próxima sprint"
$ws.Range("C18").Value2 = "Próxima sprint"
$ws.Range("C19").Value2 = "Próxima sprint"
$rng = $ws.Range("C17:C19")
$rng.Font.Italic = $true
$rng.Font.Color = 255   # pure red (BGR-packed 0x0000FF)

# Stray underline formatting that ended up on an empty cell next to the
# "Sabrine - OK" row (artifact of editing in the original workbook).
$ws.Range("F15").Font.Underline = 2

# ---------------------------------------------------------------------
# 3) SPRINT 3 content
# ---------------------------------------------------------------------

$ws.Range("A21").Value2 = "SPRINT 3"
$ws.Range("B22").Value2 = "Finalização do Front-end + Back (php)"
$ws.Range("B23").ClearContents()
$ws.Range("B24").Value2 = "Em andamento"

$ws.Range("B26").Value2 = "IMAGENS: Seleção de imagens para todas as páginas"
$ws.Range("C26").Value2 = "Stela"
$ws.Range("B27").Value2 = "TEXTOS E TÍTULOS: Definição dos textos para todas as páginas"
$ws.Range("C27").Value2 = "Stela"
$ws.Range("B28").Value2 = "LOGO: Elaboração do Logo"
$ws.Range("C28").Value2 = "Sabrine"
$ws.Range("B29").Value2 = "Ajustar páginas com modelo da YELLOW"
$ws.Range("C29").ClearContents()
$ws.Range("B30").Value2 = "Melhorar CSS (cores, fontes, etc..)"
$ws.Range("C30").ClearContents()
$ws.Range("B31").Value2 = "Login - virar modal (php) - idem modelo dado em aula"
$ws.Range("C31").ClearContents()
$ws.Range("B32").Value2 = "Criar página de alteração do cadastro - após login"
$ws.Range("C32").ClearContents()
$ws.Range("B33").Value2 = "Separar HEAD, HEADER e FOOTER - e incluir nas páginas com php"
$ws.Range("C33").ClearContents()
$ws.Range("B34").Value2 = "PHP - tudo, nem sei por onde começar a listar isso......."
$ws.Range("C34").ClearContents()

# Column C got a little wider once it started holding people's names.
$ws.Columns.Item(3).AutoFit()

# Scroll the view down to where the action now is.
$ws.Range("F31").Select()
